$wb = $excel.ActiveWorkbook

# ALC!row2 - Mercury Rising
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 331.93332
$ws.Range("I2").Value = 190
$ws.Range("J2").Value = 615.8
$ws.Range("K2").Value = 190
$ws.Range("L2").Value = 615.8
$ws.Range("M2").Value = -77
$ws.Range("N2").Value = -841.8

# ALC!row9 - Distill, My Heart
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 174.875
$ws.Range("I9").Value = 137.25
$ws.Range("J9").Value = 212.5
$ws.Range("K9").Value = 137.25
$ws.Range("L9").Value = 212.5
$ws.Range("M9").Value = 31.75
$ws.Range("N9").Value = -550.5

# ALC!row28 - The Writing Is Not on the Wall
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1638.5
$ws.Range("I28").Value = 1446.7646
$ws.Range("J28").Value = 2290.4
$ws.Range("K28").Value = 1446.7646
$ws.Range("L28").Value = 2290.4
$ws.Range("M28").Value = -961.7646
$ws.Range("N28").Value = -3260.4

# ALC!row58 - A Matter of Vital Importance
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1001.1905
$ws.Range("I58").Value = 817.17645
$ws.Range("J58").Value = 1783.25
$ws.Range("K58").Value = 2451.52935
$ws.Range("L58").Value = 5349.75
$ws.Range("M58").Value = -2301.52935
$ws.Range("N58").Value = -5649.75

# ALC!row100 - Asking for a Friend
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2440.8333
$ws.Range("I100").Value = 1869.1428
$ws.Range("J100").Value = 3241.2
$ws.Range("K100").Value = 1869.1428
$ws.Range("L100").Value = 3241.2
$ws.Range("M100").Value = -1328.1428
$ws.Range("N100").Value = -4323.2

# ALC!row107 - Another Man's Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 459.7647
$ws.Range("I107").Value = 480
$ws.Range("J107").Value = 136
$ws.Range("K107").Value = 480
$ws.Range("L107").Value = 136
$ws.Range("M107").Value = 1440
$ws.Range("N107").Value = -3976

# ALC!row125 - Body over Mind
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 979.38464
$ws.Range("I125").Value = 842.4286
$ws.Range("J125").Value = 1139.1666
$ws.Range("K125").Value = 7581.8574
$ws.Range("L125").Value = 10252.4994
$ws.Range("M125").Value = -5121.8574
$ws.Range("N125").Value = -15172.4994

# ALC!row132 - Fast-forwarding Flora
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3306.513
$ws.Range("I132").Value = 1712.7858
$ws.Range("K132").Value = 5138.357400000001
$ws.Range("M132").Value = -2608.357400000001

# ARM!row11 - Rodents of Unusual Size
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 2500250
$ws.Range("J11").Value = 500
$ws.Range("L11").Value = 500
$ws.Range("N11").Value = -788

# ARM!row32 - Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2853.2705
$ws.Range("I32").Value = 2272.5063
$ws.Range("K32").Value = 2272.5063
$ws.Range("M32").Value = -1985.5063

# BSM!row20 - Smelt and Dealt
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1690.8286
$ws.Range("I20").Value = 1337.625
$ws.Range("J20").Value = 1988.2632
$ws.Range("K20").Value = 1337.625
$ws.Range("L20").Value = 1988.2632
$ws.Range("M20").Value = -1090.625
$ws.Range("N20").Value = -2482.2632

# BSM!row94 - High Steal
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 597.25
$ws.Range("I94").Value = 519.61536
$ws.Range("J94").Value = 741.4286
$ws.Range("K94").Value = 519.61536
$ws.Range("L94").Value = 741.4286
$ws.Range("M94").Value = -68.61536000000001
$ws.Range("N94").Value = -1643.4286

# BSM!row126 - Records of the Republic
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 34749.125
$ws.Range("J126").Value = 34749.125
$ws.Range("L126").Value = 34749.125
$ws.Range("N126").Value = -44629.125

# CRP!row31 - Wall Not Found
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41668150
$ws.Range("I31").Value = 40001070
$ws.Range("J31").Value = 45456972
$ws.Range("K31").Value = 40001070
$ws.Range("L31").Value = 45456972
$ws.Range("M31").Value = -40000775
$ws.Range("N31").Value = -45457562

# CRP!row34 - Armoires of the Rich and Famous
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 41668150
$ws.Range("I34").Value = 40001070
$ws.Range("J34").Value = 45456972
$ws.Range("K34").Value = 40001070
$ws.Range("L34").Value = 45456972
$ws.Range("M34").Value = -40000868
$ws.Range("N34").Value = -45457376

# CRP!row107 - Built to Last
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1515.5
$ws.Range("I107").Value = 453.14285
$ws.Range("J107").Value = 1952.9412
$ws.Range("K107").Value = 453.14285
$ws.Range("L107").Value = 1952.9412
$ws.Range("M107").Value = 1466.85715
$ws.Range("N107").Value = -5792.9412

# CRP!row132 - Hull Lotta Damage
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2458.0852
$ws.Range("I132").Value = 2253
$ws.Range("J132").Value = 2895.6
$ws.Range("K132").Value = 6759
$ws.Range("L132").Value = 8686.799999999999
$ws.Range("M132").Value = -4229
$ws.Range("N132").Value = -13746.8

# CUL!row5 - What a Sap
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 728.44446
$ws.Range("I5").Value = 470.15384
$ws.Range("K5").Value = 1410.46152
$ws.Range("M5").Value = -1298.46152

# CUL!row111 - Soup for the Soldier
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 11675.667
$ws.Range("I111").Value = 12910.125
$ws.Range("J111").Value = 1800
$ws.Range("K111").Value = 38730.375
$ws.Range("L111").Value = 5400
$ws.Range("M111").Value = -35663.375
$ws.Range("N111").Value = -11534

# CUL!row135 - Not-so-secret Ingredient
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 728.44446
$ws.Range("I135").Value = 470.15384
$ws.Range("K135").Value = 4231.38456
$ws.Range("M135").Value = -1696.38456

# GSM!row70 - Sky Is the Limit
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5330.393
$ws.Range("I70").Value = 5100.5
$ws.Range("J70").Value = 5636.9165
$ws.Range("K70").Value = 5100.5
$ws.Range("L70").Value = 5636.9165
$ws.Range("M70").Value = -4830.5
$ws.Range("N70").Value = -6176.9165

# GSM!row73 - Hulls of Broken Dreams (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5330.393
$ws.Range("I73").Value = 5100.5
$ws.Range("J73").Value = 5636.9165
$ws.Range("K73").Value = 5100.5
$ws.Range("L73").Value = 5636.9165
$ws.Range("M73").Value = -4164.5
$ws.Range("N73").Value = -7508.9165

# GSM!row107 - Whetstones for the Workers
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 754.38464
$ws.Range("J107").Value = 801.7692
$ws.Range("L107").Value = 801.7692
$ws.Range("N107").Value = -4641.7692

# GSM!row113 - Copious Crystal Cannons
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2037.1111
$ws.Range("I113").Value = 1249.5
$ws.Range("J113").Value = 2262.1428
$ws.Range("K113").Value = 1249.5
$ws.Range("L113").Value = 2262.1428
$ws.Range("M113").Value = 920.5
$ws.Range("N113").Value = -6602.1428

# LTW!row16 - Saddle Sore
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 671.4666999999999
$ws.Range("I16").Value = 690.7143
$ws.Range("K16").Value = 690.7143
$ws.Range("M16").Value = -520.7143

# LTW!row61 - Spelling Me Softly
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2601.6
$ws.Range("I61").Value = 1669.3334
$ws.Range("K61").Value = 1669.3334
$ws.Range("M61").Value = -1467.3334

# LTW!row113 - Peace in Rest
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2601.6
$ws.Range("I113").Value = 1669.3334
$ws.Range("K113").Value = 1669.3334
$ws.Range("M113").Value = 500.6666

# WVR!row132 - Comfy Cabins
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8338062.5
$ws.Range("I132").Value = 12505231
$ws.Range("J132").Value = 3724.875
$ws.Range("K132").Value = 37515693
$ws.Range("L132").Value = 11174.625
$ws.Range("M132").Value = -37513163
$ws.Range("N132").Value = -16234.625
